{"js": "// Update the date heading and the 25 division-fact answers in the table.\n// Each edit replaces only the text of an existing paragraph, leaving the\n// run/paragraph formatting (fonts, size, alignment) untouched.\n\n// 1) Date heading paragraph, e.g. \"2024-04-08 Monday\" -> \"2024-04-09 Tuesday\"\nconst headingParas = context.document.body.paragraphs;\nheadingParas.load(\"items\");\nawait context.sync();\n\nconst heading = headingParas.items[0];\nheading.insertText(\"2024-04-09 Tuesday\", Word.InsertLocation.replace);\n\n// 2) The division problems live in a single 5-column table; only every 4th\n// row (0, 4, 8, 12, 16) actually holds text, the rest are spacer rows.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New answers, in row-major order matching the rows that contain text.\nconst newValues = [\n  [\"16\u00f77=2, 2\", \"88\u00f74=22, 0\", \"81\u00f79=9, 0\", \"79\u00f79=8, 7\", \"39\u00f77=5, 4\"],\n  [\"91\u00f73=30, 1\", \"66\u00f76=11, 0\", \"34\u00f76=5, 4\", \"44\u00f76=7, 2\", \"39\u00f73=13, 0\"],\n  [\"19\u00f72=9, 1\", \"59\u00f78=7, 3\", \"40\u00f73=13, 1\", \"82\u00f76=13, 4\", \"93\u00f75=18, 3\"],\n  [\"51\u00f78=6, 3\", \"57\u00f77=8, 1\", \"54\u00f76=9, 0\", \"21\u00f77=3, 0\", \"23\u00f76=3, 5\"],\n  [\"79\u00f74=19, 3\", \"52\u00f79=5, 7\", \"56\u00f76=9, 2\", \"69\u00f78=8, 5\", \"60\u00f77=8, 4\"],\n];\nconst contentRowIndexes = [0, 4, 8, 12, 16];\n\nfor (let r = 0; r < contentRowIndexes.length; r++) {\n  const rowIndex = contentRowIndexes[r];\n  for (let c = 0; c < 5; c++) {\n    const cell = table.getCell(rowIndex, c);\n    const cellParas = cell.body.paragraphs;\n    cellParas.load(\"items\");\n    await context.sync();\n    const p = cellParas.items[0];\n    p.insertText(newValues[r][c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 division-fact answers in the table.\n# Each assignment replaces only the run text, leaving paragraph/run\n# formatting (fonts, size, alignment) untouched.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading paragraph, e.g. \"2024-04-08 Monday\" -> \"2024-04-09 Tuesday\"\n$d.Paragraphs.Item(1).Range.Text = \"2024-04-09 Tuesday\"\n\n# 2) The division problems live in a single 5-column table; only every 4th\n# row (1, 5, 9, 13, 17 in 1-based Word numbering) actually holds text, the\n# rest are spacer rows.\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"16\u00f77=2, 2\", \"88\u00f74=22, 0\", \"81\u00f79=9, 0\", \"79\u00f79=8, 7\", \"39\u00f77=5, 4\"),\n    @(\"91\u00f73=30, 1\", \"66\u00f76=11, 0\", \"34\u00f76=5, 4\", \"44\u00f76=7, 2\", \"39\u00f73=13, 0\"),\n    @(\"19\u00f72=9, 1\", \"59\u00f78=7, 3\", \"40\u00f73=13, 1\", \"82\u00f76=13, 4\", \"93\u00f75=18, 3\"),\n    @(\"51\u00f78=6, 3\", \"57\u00f77=8, 1\", \"54\u00f76=9, 0\", \"21\u00f77=3, 0\", \"23\u00f76=3, 5\"),\n    @(\"79\u00f74=19, 3\", \"52\u00f79=5, 7\", \"56\u00f76=9, 2\", \"69\u00f78=8, 5\", \"60\u00f77=8, 4\")\n)\n$contentRows = @(1, 5, 9, 13, 17)\n\nfor ($r = 0; $r -lt $contentRows.Length; $r++) {\n    $rowIndex = $contentRows[$r]\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $t.Cell($rowIndex, $c)\n        $cell.Range.Text = $newValues[$r][$c - 1]\n    }\n}\n"}
